$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 119285.766
$ws.Cells.Item(17, 10).Value = 119285.766
$ws.Cells.Item(17, 12).Value = 357857.298
$ws.Cells.Item(17, 14).Value = -358193.298
$ws.Cells.Item(112, 8).Value = 43884.51
$ws.Cells.Item(112, 10).Value = 30467.568
$ws.Cells.Item(112, 12).Value = 91402.704
$ws.Cells.Item(112, 14).Value = -93618.704
$ws.Cells.Item(125, 8).Value = 13892943
$ws.Cells.Item(125, 9).Value = 4289.8335
$ws.Cells.Item(125, 11).Value = 38608.5015
$ws.Cells.Item(125, 13).Value = -36148.5015
$ws.Cells.Item(132, 8).Value = 3410.9524
$ws.Cells.Item(132, 9).Value = 3975.3845
$ws.Cells.Item(132, 11).Value = 11926.1535
$ws.Cells.Item(132, 13).Value = -9396.1535
$ws.Cells.Item(137, 8).Value = 1309.5416
$ws.Cells.Item(137, 10).Value = 3800.75
$ws.Cells.Item(137, 12).Value = 11402.25
$ws.Cells.Item(137, 14).Value = -16502.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 866937.0600000001
$ws.Cells.Item(2, 9).Value = 1132719
$ws.Cells.Item(2, 10).Value = 3145.75
$ws.Cells.Item(2, 11).Value = 1132719
$ws.Cells.Item(2, 12).Value = 3145.75
$ws.Cells.Item(2, 13).Value = -1132606
$ws.Cells.Item(2, 14).Value = -3371.75
$ws.Cells.Item(61, 8).Value = 52634760
$ws.Cells.Item(61, 9).Value = 62503056
$ws.Cells.Item(61, 11).Value = 62503056
$ws.Cells.Item(61, 13).Value = -62502844
$ws.Cells.Item(74, 8).Value = 76930790
$ws.Cells.Item(74, 9).Value = 76930790
$ws.Cells.Item(74, 11).Value = 76930790
$ws.Cells.Item(74, 13).Value = -76929916
$ws.Cells.Item(77, 8).Value = 76930790
$ws.Cells.Item(77, 9).Value = 76930790
$ws.Cells.Item(77, 11).Value = 384653950
$ws.Cells.Item(77, 13).Value = -384649582
$ws.Cells.Item(116, 8).Value = 866937.0600000001
$ws.Cells.Item(116, 9).Value = 1132719
$ws.Cells.Item(116, 10).Value = 3145.75
$ws.Cells.Item(116, 11).Value = 1132719
$ws.Cells.Item(116, 12).Value = 3145.75
$ws.Cells.Item(116, 13).Value = -1130425
$ws.Cells.Item(116, 14).Value = -7733.75
$ws.Cells.Item(132, 8).Value = 2566648.5
$ws.Cells.Item(132, 9).Value = 2859354.2
$ws.Cells.Item(132, 10).Value = 5474
$ws.Cells.Item(132, 11).Value = 8578062.600000001
$ws.Cells.Item(132, 12).Value = 16422
$ws.Cells.Item(132, 13).Value = -8575532.600000001
$ws.Cells.Item(132, 14).Value = -21482
$ws.Cells.Item(136, 8).Value = 52634760
$ws.Cells.Item(136, 9).Value = 62503056
$ws.Cells.Item(136, 11).Value = 187509168
$ws.Cells.Item(136, 13).Value = -187506618

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 866937.0600000001
$ws.Cells.Item(3, 9).Value = 1132719
$ws.Cells.Item(3, 10).Value = 3145.75
$ws.Cells.Item(3, 11).Value = 1132719
$ws.Cells.Item(3, 12).Value = 3145.75
$ws.Cells.Item(3, 13).Value = -1132605
$ws.Cells.Item(3, 14).Value = -3373.75
$ws.Cells.Item(19, 8).Value = 20000
$ws.Cells.Item(19, 10).Value = 20000
$ws.Cells.Item(19, 12).Value = 20000
$ws.Cells.Item(19, 14).Value = -20346
$ws.Cells.Item(94, 8).Value = 938.36365
$ws.Cells.Item(94, 9).Value = 942.1667
$ws.Cells.Item(94, 11).Value = 942.1667
$ws.Cells.Item(94, 13).Value = -491.1667
$ws.Cells.Item(105, 8).Value = 1783.0952
$ws.Cells.Item(105, 9).Value = 1748.1111
$ws.Cells.Item(105, 11).Value = 1748.1111
$ws.Cells.Item(105, 13).Value = -1.111100000000079
$ws.Cells.Item(134, 8).Value = 40387664
$ws.Cells.Item(134, 9).Value = 47730084
$ws.Cells.Item(134, 10).Value = 4347
$ws.Cells.Item(134, 11).Value = 143190252
$ws.Cells.Item(134, 12).Value = 13041
$ws.Cells.Item(134, 13).Value = -143187717
$ws.Cells.Item(134, 14).Value = -18111

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 362.5
$ws.Cells.Item(2, 9).Value = 350
$ws.Cells.Item(2, 11).Value = 350
$ws.Cells.Item(2, 13).Value = -237
$ws.Cells.Item(16, 8).Value = 1089230.2
$ws.Cells.Item(16, 9).Value = 1812200.6
$ws.Cells.Item(16, 11).Value = 1812200.6
$ws.Cells.Item(16, 13).Value = -1811913.6
$ws.Cells.Item(58, 8).Value = 62514864
$ws.Cells.Item(58, 9).Value = 71445310
$ws.Cells.Item(58, 11).Value = 71445310
$ws.Cells.Item(58, 13).Value = -71445107
$ws.Cells.Item(99, 8).Value = 1914.5
$ws.Cells.Item(99, 9).Value = 1798.4
$ws.Cells.Item(99, 10).Value = 2495
$ws.Cells.Item(99, 11).Value = 1798.4
$ws.Cells.Item(99, 12).Value = 2495
$ws.Cells.Item(99, 13).Value = -300.4000000000001
$ws.Cells.Item(99, 14).Value = -5491
$ws.Cells.Item(113, 8).Value = 1089230.2
$ws.Cells.Item(113, 9).Value = 1812200.6
$ws.Cells.Item(113, 11).Value = 1812200.6
$ws.Cells.Item(113, 13).Value = -1810030.6
$ws.Cells.Item(126, 8).Value = 1914.5
$ws.Cells.Item(126, 9).Value = 1798.4
$ws.Cells.Item(126, 10).Value = 2495
$ws.Cells.Item(126, 11).Value = 5395.200000000001
$ws.Cells.Item(126, 12).Value = 7485
$ws.Cells.Item(126, 13).Value = -2925.200000000001
$ws.Cells.Item(126, 14).Value = -12425
$ws.Cells.Item(132, 8).Value = 125002880
$ws.Cells.Item(132, 9).Value = 166670060
$ws.Cells.Item(132, 11).Value = 500010180
$ws.Cells.Item(132, 13).Value = -500007650
$ws.Cells.Item(134, 8).Value = 31252318
$ws.Cells.Item(134, 9).Value = 41668424
$ws.Cells.Item(134, 11).Value = 125005272
$ws.Cells.Item(134, 13).Value = -125002737
$ws.Cells.Item(136, 8).Value = 62514864
$ws.Cells.Item(136, 9).Value = 71445310
$ws.Cells.Item(136, 11).Value = 214335930
$ws.Cells.Item(136, 13).Value = -214333380
$ws.Cells.Item(140, 8).Value = 86239.75
$ws.Cells.Item(140, 10).Value = 86239.75
$ws.Cells.Item(140, 12).Value = 86239.75
$ws.Cells.Item(140, 14).Value = -96599.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(117, 8).Value = 5649.5
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 10).Value = 5649.5
$ws.Cells.Item(117, 11).Value = 0
$ws.Cells.Item(117, 12).Value = 16948.5
$ws.Cells.Item(117, 13).ClearContents()
$ws.Cells.Item(117, 14).Value = -23832.5
$ws.Cells.Item(140, 8).Value = 1757.1923
$ws.Cells.Item(140, 9).Value = 722.2273
$ws.Cells.Item(140, 11).Value = 2166.6819
$ws.Cells.Item(140, 13).Value = 3013.3181
$ws.Cells.Item(141, 8).Value = 1196
$ws.Cells.Item(141, 9).Value = 1196
$ws.Cells.Item(141, 11).Value = 3588
$ws.Cells.Item(141, 13).Value = 1592

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 508.1613
$ws.Cells.Item(97, 9).Value = 349.42105
$ws.Cells.Item(97, 10).Value = 759.5
$ws.Cells.Item(97, 11).Value = 349.42105
$ws.Cells.Item(97, 12).Value = 759.5
$ws.Cells.Item(97, 13).Value = 146.57895
$ws.Cells.Item(97, 14).Value = -1751.5
$ws.Cells.Item(100, 8).Value = 130000
$ws.Cells.Item(100, 10).Value = 130000
$ws.Cells.Item(100, 12).Value = 130000
$ws.Cells.Item(100, 14).Value = -132164
$ws.Cells.Item(102, 8).Value = 1964.1111
$ws.Cells.Item(102, 9).Value = 1783.1428
$ws.Cells.Item(102, 10).Value = 2597.5
$ws.Cells.Item(102, 11).Value = 1783.1428
$ws.Cells.Item(102, 12).Value = 2597.5
$ws.Cells.Item(102, 13).Value = -161.1428000000001
$ws.Cells.Item(102, 14).Value = -5841.5
$ws.Cells.Item(132, 8).Value = 6946976
$ws.Cells.Item(132, 9).Value = 7355605.5
$ws.Cells.Item(132, 11).Value = 22066816.5
$ws.Cells.Item(132, 13).Value = -22064286.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2824.5
$ws.Cells.Item(82, 9).Value = 2824.5
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 2824.5
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 13).Value = -2463.5
$ws.Cells.Item(82, 14).ClearContents()
$ws.Cells.Item(85, 8).Value = 2824.5
$ws.Cells.Item(85, 9).Value = 2824.5
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 2824.5
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = -1576.5
$ws.Cells.Item(85, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 17866684
$ws.Cells.Item(132, 9).Value = 17866684
$ws.Cells.Item(132, 11).Value = 53600052
$ws.Cells.Item(132, 13).Value = -53597522
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 14).ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 6830.3335
$ws.Cells.Item(122, 9).Value = 8003.6665
$ws.Cells.Item(122, 10).Value = 3897
$ws.Cells.Item(122, 11).Value = 24010.9995
$ws.Cells.Item(122, 12).Value = 11691
$ws.Cells.Item(122, 13).Value = -21560.9995
$ws.Cells.Item(122, 14).Value = -16591
$ws.Cells.Item(132, 8).Value = 15630496
$ws.Cells.Item(132, 9).Value = 20835530
$ws.Cells.Item(132, 11).Value = 62506590
$ws.Cells.Item(132, 13).Value = -62504060
$ws.Cells.Item(136, 8).Value = 12502670
$ws.Cells.Item(136, 10).Value = 3263.6667
$ws.Cells.Item(136, 12).Value = 9791.000100000001
$ws.Cells.Item(136, 14).Value = -14891.0001
